# Auto-generated script applying numeric cell updates described in the commit diff.
# Updates LeveProfit/price/average-price columns (H-N) across several worksheets.

$wb = $excel.ActiveWorkbook

# --- ALC!row33 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 290.15
$ws.Range("J33").Value = 663.3333
$ws.Range("L33").Value = 663.3333
$ws.Range("N33").Value = -1121.3333

# --- ALC!row38 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 845.625
$ws.Range("J38").Value = 425
$ws.Range("L38").Value = 1275
$ws.Range("N38").Value = -2019

# --- ALC!row40 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3326.6858
$ws.Range("I40").Value = 1818.8572
$ws.Range("J40").Value = 3703.6428
$ws.Range("K40").Value = 1818.8572
$ws.Range("L40").Value = 3703.6428
$ws.Range("M40").Value = -1643.8572
$ws.Range("N40").Value = -4053.6428

# --- ALC!row58 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 3018.4285
$ws.Range("I58").Value = 2682.25
$ws.Range("J58").Value = 3466.6667
$ws.Range("K58").Value = 8046.75
$ws.Range("L58").Value = 10400.0001
$ws.Range("M58").Value = -7896.75
$ws.Range("N58").Value = -10700.0001

# --- ALC!row61 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 269.57144
$ws.Range("I61").Value = 269.57144
$ws.Range("K61").Value = 808.71432
$ws.Range("M61").Value = -636.71432

# --- ALC!row74 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5216.467
$ws.Range("I74").Value = 5219
$ws.Range("K74").Value = 5219
$ws.Range("M74").Value = -4283

# --- ALC!row76 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 19348.2
$ws.Range("I76").Value = 29397.6
$ws.Range("K76").Value = 29397.6
$ws.Range("M76").Value = -29082.6

# --- ALC!row77 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 5216.467
$ws.Range("I77").Value = 5219
$ws.Range("K77").Value = 26095
$ws.Range("M77").Value = -21415

# --- ALC!row79 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 19348.2
$ws.Range("I79").Value = 29397.6
$ws.Range("K79").Value = 29397.6
$ws.Range("M79").Value = -28305.6

# --- ALC!row98 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2440.4
$ws.Range("I98").Value = 1343.5714
$ws.Range("K98").Value = 1343.5714
$ws.Range("M98").Value = 154.4286

# --- ALC!row122 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2440.4
$ws.Range("I122").Value = 1343.5714
$ws.Range("K122").Value = 4030.7142
$ws.Range("M122").Value = -1580.7142

# --- ALC!row132 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1811.7407
$ws.Range("I132").Value = 1813.2916
$ws.Range("K132").Value = 5439.8748
$ws.Range("M132").Value = -2909.8748

# --- ALC!row135 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2820.182
$ws.Range("I135").Value = 1802.2
$ws.Range("K135").Value = 16219.8
$ws.Range("M135").Value = -13684.8

# --- ARM!row122 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3589.7
$ws.Range("I122").Value = 3245.9092
$ws.Range("J122").Value = 4009.889
$ws.Range("K122").Value = 9737.7276
$ws.Range("L122").Value = 12029.667
$ws.Range("M122").Value = -7287.7276
$ws.Range("N122").Value = -16929.667

# --- BSM!row22 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -127

# --- BSM!row134 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1590189.6
$ws.Range("I134").Value = 1985926
$ws.Range("J134").Value = 7243.6665
$ws.Range("K134").Value = 5957778
$ws.Range("L134").Value = 21730.9995
$ws.Range("M134").Value = -5955243
$ws.Range("N134").Value = -26800.9995

# --- CRP!row58 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2753.149
$ws.Range("I58").Value = 2521.3
$ws.Range("K58").Value = 2521.3
$ws.Range("M58").Value = -2318.3

# --- CRP!row86 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 34014.89
$ws.Range("I86").Value = 24008.904
$ws.Range("K86").Value = 24008.904
$ws.Range("M86").Value = -22885.904

# --- CRP!row89 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 34014.89
$ws.Range("I89").Value = 24008.904
$ws.Range("K89").Value = 120044.52
$ws.Range("M89").Value = -114428.52

# --- CRP!row132 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3091.4473
$ws.Range("I132").Value = 3058.0334
$ws.Range("J132").Value = 3216.75
$ws.Range("K132").Value = 9174.100199999999
$ws.Range("L132").Value = 9650.25
$ws.Range("M132").Value = -6644.100199999999
$ws.Range("N132").Value = -14710.25

# --- CRP!row134 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3704.647
$ws.Range("I134").Value = 3704.647
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 11113.941
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -8578.940999999999
$ws.Range("N134").ClearContents()

# --- CRP!row136 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2753.149
$ws.Range("I136").Value = 2521.3
$ws.Range("K136").Value = 7563.900000000001
$ws.Range("M136").Value = -5013.900000000001

# --- CUL!row18 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 4934.1113
$ws.Range("I18").Value = 4081.6
$ws.Range("K18").Value = 12244.8
$ws.Range("M18").Value = -12075.8

# --- CUL!row132 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2507299
$ws.Range("I132").Value = 9200
$ws.Range("K132").Value = 82800
$ws.Range("M132").Value = -80270

# --- GSM!row70 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 24196.75
$ws.Range("I70").Value = 66329.7
$ws.Range("K70").Value = 66329.7
$ws.Range("M70").Value = -66059.7

# --- GSM!row73 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 24196.75
$ws.Range("I73").Value = 66329.7
$ws.Range("K73").Value = 66329.7
$ws.Range("M73").Value = -65393.7

# --- GSM!row102 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3286.7778
$ws.Range("I102").Value = 3197.625
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 3197.625
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -1575.625
$ws.Range("N102").Value = -7244

# --- GSM!row122 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1044.6
$ws.Range("I122").Value = 1044.6
$ws.Range("K122").Value = 3133.8
$ws.Range("M122").Value = -683.7999999999997

# --- GSM!row126 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2667.6667
$ws.Range("I126").Value = 2501.2
$ws.Range("K126").Value = 7503.599999999999
$ws.Range("M126").Value = -5033.599999999999

# --- LTW!row68 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8146.3076
$ws.Range("I68").Value = 6933.222
$ws.Range("J68").Value = 10875.75
$ws.Range("K68").Value = 6933.222
$ws.Range("L68").Value = 10875.75
$ws.Range("M68").Value = -6184.222
$ws.Range("N68").Value = -12373.75

# --- LTW!row71 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 8146.3076
$ws.Range("I71").Value = 6933.222
$ws.Range("J71").Value = 10875.75
$ws.Range("K71").Value = 34666.11
$ws.Range("L71").Value = 54378.75
$ws.Range("M71").Value = -30922.11
$ws.Range("N71").Value = -61866.75

# --- LTW!row82 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4179.385
$ws.Range("J82").Value = 3838.8333
$ws.Range("L82").Value = 3838.8333
$ws.Range("N82").Value = -4560.8333

# --- LTW!row85 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 4179.385
$ws.Range("J85").Value = 3838.8333
$ws.Range("L85").Value = 3838.8333
$ws.Range("N85").Value = -6334.8333

# --- LTW!row122 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 21179.53
$ws.Range("J122").Value = 24000.5
$ws.Range("L122").Value = 72001.5
$ws.Range("N122").Value = -76901.5

# --- LTW!row132 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2865.65
$ws.Range("I132").Value = 2753.3157
$ws.Range("K132").Value = 8259.947100000001
$ws.Range("M132").Value = -5729.947100000001

# --- LTW!row136 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1831.421
$ws.Range("I136").Value = 1699.5
$ws.Range("J136").Value = 2535
$ws.Range("K136").Value = 5098.5
$ws.Range("L136").Value = 7605
$ws.Range("M136").Value = -2548.5
